# Apply the "Final epsilons" value updates to the overall_scores worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2"  = 0.5649999999999999
    "E2"  = 0.293
    "H2"  = 0.525
    "K2"  = 0.707
    "L2"  = 1
    "N2"  = 0.716
    "Q2"  = 0

    "B3"  = 0.828
    "E3"  = 0.317
    "H3"  = 0.437
    "K3"  = 0.853
    "L3"  = 1
    "N3"  = 0.763

    "B4"  = 0.5649999999999999
    "E4"  = 0.165
    "H4"  = 0.323
    "K4"  = 0.67
    "N4"  = 0.739

    "B5"  = 0.9409999999999999
    "E5"  = 0.796
    "H5"  = 0.854
    "K5"  = 0.978
    "N5"  = 0.235

    "B6"  = 0.971
    "E6"  = 0.965
    "H6"  = 0.328
    "K6"  = 0.99
    "N6"  = 0.2

    "B7"  = 0.877
    "E7"  = 0.803
    "H7"  = 0.591
    "K7"  = 0.867
    "N7"  = 0.502

    "B8"  = 0.91
    "E8"  = 0.841
    "H8"  = 0.975
    "K8"  = 0.993
    "N8"  = 0.845

    "B9"  = 0.205
    "E9"  = 0.908
    "H9"  = 0.861
    "K9"  = 0.994
    "N9"  = 0.728

    "B10" = 0.5649999999999999
    "E10" = 0.876
    "H10" = 0.953
    "K10" = 0.9370000000000001
    "N10" = 0.694

    "B11" = 0.902
    "E11" = 0.821
    "H11" = 0.856
    "K11" = 0.99
    "N11" = 0.457
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
